# Populate Sheet1 with the new cell values.
# Shared-string table indices in the target workbook are assigned in the
# order the strings are first written (fdgsdf=0, sdfgfd=1, sdgfdg=2,
# sdgfd=3, asdfdsas=4), so the writes below follow that exact order.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = "fdgsdf"
$ws.Range("D7").Value = "sdfgfd"
$ws.Range("F10").Value = "sdgfdg"
$ws.Range("D3").Value = "sdgfd"
$ws.Range("K10").Value = "asdfdsas"

# Final selection lands on K10, matching the saved view state.
$ws.Range("K10").Select() | Out-Null
